$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -18.07964102843893
$ws.Range("C2").Value = 1.9435948822869
$ws.Range("D2").Value = -18.07964102843893
$ws.Range("E2").Value = -18.07964102843893
$ws.Range("F2").Value = -18.07964102843893
$ws.Range("G2").Value = -18.07964102843893
$ws.Range("H2").Value = -18.07964102843893
$ws.Range("I2").Value = -18.07964102843893
$ws.Range("J2").Value = -18.07964102843893
$ws.Range("K2").Value = -18.07964102843893
$ws.Range("B3").Value = -18.07964102843893
$ws.Range("C3").Value = -18.07964102843893
$ws.Range("D3").Value = -18.07964102843893
$ws.Range("E3").Value = -18.07964102843893
$ws.Range("F3").Value = -18.07964102843893
$ws.Range("G3").Value = -18.07964102843893
$ws.Range("H3").Value = -18.07964102843893
$ws.Range("I3").Value = -0.4300129400115966
$ws.Range("J3").Value = -18.07964102843893
$ws.Range("K3").Value = -18.07964102843893
$ws.Range("B4").Value = -18.07964102843893
$ws.Range("C4").Value = 1.937326361651383
$ws.Range("D4").Value = 1.562207869991699
$ws.Range("E4").Value = -18.07964102843893
$ws.Range("F4").Value = 3.252105298228955
$ws.Range("G4").Value = -18.07964102843893
$ws.Range("H4").Value = 1.328508314408298
$ws.Range("I4").Value = -18.07964102843893
$ws.Range("J4").Value = 4.32192314740074
$ws.Range("K4").Value = -18.07964102843893
$ws.Range("B5").Value = -18.07964102843893
$ws.Range("C5").Value = 1.728149386230838
$ws.Range("D5").Value = -18.07964102843893
$ws.Range("E5").Value = -18.07964102843893
$ws.Range("F5").Value = -18.07964102843893
$ws.Range("G5").Value = 2.656014977173438
$ws.Range("H5").Value = -18.07964102843893
$ws.Range("I5").Value = -18.07964102843893
$ws.Range("J5").Value = -18.07964102843893
$ws.Range("K5").Value = -18.07964102843893
$ws.Range("B6").Value = -18.07964102843893
$ws.Range("C6").Value = -18.07964102843893
$ws.Range("D6").Value = -18.07964102843893
$ws.Range("E6").Value = -18.07964102843893
$ws.Range("F6").Value = -18.07964102843893
$ws.Range("G6").Value = -18.07964102843893
$ws.Range("H6").Value = -18.07964102843893
$ws.Range("I6").Value = -18.07964102843893
$ws.Range("J6").Value = -18.07964102843893
$ws.Range("K6").Value = -18.07964102843893
$ws.Range("B7").Value = 2.438294542695173
$ws.Range("C7").Value = -18.07964102843893
$ws.Range("D7").Value = -18.07964102843893
$ws.Range("E7").Value = -18.07964102843893
$ws.Range("F7").Value = -18.07964102843893
$ws.Range("G7").Value = -18.07964102843893
$ws.Range("H7").Value = -18.07964102843893
$ws.Range("I7").Value = -18.07964102843893
$ws.Range("J7").Value = -18.07964102843893
$ws.Range("K7").Value = -18.07964102843893
$ws.Range("B8").Value = -18.07964102843893
$ws.Range("C8").Value = -18.07964102843893
$ws.Range("D8").Value = -18.07964102843893
$ws.Range("E8").Value = 1.726608063712083
$ws.Range("F8").Value = -18.07964102843893
$ws.Range("G8").Value = -18.07964102843893
$ws.Range("H8").Value = -18.07964102843893
$ws.Range("I8").Value = -18.07964102843893
$ws.Range("J8").Value = -18.07964102843893
$ws.Range("K8").Value = -18.07964102843893
$ws.Range("B9").Value = 3.865911756404269
$ws.Range("C9").Value = -18.07964102843893
$ws.Range("D9").Value = -18.07964102843893
$ws.Range("E9").Value = -18.07964102843893
$ws.Range("F9").Value = -18.07964102843893
$ws.Range("G9").Value = -18.07964102843893
$ws.Range("H9").Value = -18.07964102843893
$ws.Range("I9").Value = -18.07964102843893
$ws.Range("J9").Value = -18.07964102843893
$ws.Range("K9").Value = -18.07964102843893
$ws.Range("B10").Value = -18.07964102843893
$ws.Range("C10").Value = -18.07964102843893
$ws.Range("D10").Value = -18.07964102843893
$ws.Range("E10").Value = -18.07964102843893
$ws.Range("F10").Value = -18.07964102843893
$ws.Range("G10").Value = -18.07964102843893
$ws.Range("H10").Value = -18.07964102843893
$ws.Range("I10").Value = 2.021340775727742
$ws.Range("J10").Value = -18.07964102843893
$ws.Range("K10").Value = 2.136869207580869
$ws.Range("B11").Value = -18.07964102843893
$ws.Range("C11").Value = -18.07964102843893
$ws.Range("D11").Value = -18.07964102843893
$ws.Range("E11").Value = 2.932504029034326
$ws.Range("F11").Value = -18.07964102843893
$ws.Range("G11").Value = 2.688498864172611
$ws.Range("H11").Value = -18.07964102843893
$ws.Range("I11").Value = -18.07964102843893
$ws.Range("J11").Value = -18.07964102843893
$ws.Range("K11").Value = 2.07853594612418
$ws.Range("B12").Value = -18.07964102843893
$ws.Range("C12").Value = -18.07964102843893
$ws.Range("D12").Value = -18.07964102843893
$ws.Range("E12").Value = -18.07964102843893
$ws.Range("F12").Value = -18.07964102843893
$ws.Range("G12").Value = -18.07964102843893
$ws.Range("H12").Value = -18.07964102843893
$ws.Range("I12").Value = -18.07964102843893
$ws.Range("J12").Value = -18.07964102843893
$ws.Range("K12").Value = -18.07964102843893
$ws.Range("B13").Value = -18.07964102843893
$ws.Range("C13").Value = -18.07964102843893
$ws.Range("D13").Value = -18.07964102843893
$ws.Range("E13").Value = 2.549666671301687
$ws.Range("F13").Value = -18.07964102843893
$ws.Range("G13").Value = -18.07964102843893
$ws.Range("H13").Value = -18.07964102843893
$ws.Range("I13").Value = -18.07964102843893
$ws.Range("J13").Value = -18.07964102843893
$ws.Range("K13").Value = 1.178697266161396
$ws.Range("B14").Value = -18.07964102843893
$ws.Range("C14").Value = -18.07964102843893
$ws.Range("D14").Value = 1.509002246426331
$ws.Range("E14").Value = -18.07964102843893
$ws.Range("F14").Value = -18.07964102843893
$ws.Range("G14").Value = -18.07964102843893
$ws.Range("H14").Value = -18.07964102843893
$ws.Range("I14").Value = -18.07964102843893
$ws.Range("J14").Value = -18.07964102843893
$ws.Range("K14").Value = 1.846552156293417
$ws.Range("B15").Value = -18.07964102843893
$ws.Range("C15").Value = -18.07964102843893
$ws.Range("D15").Value = 1.794387246230622
$ws.Range("E15").Value = -18.07964102843893
$ws.Range("F15").Value = -18.07964102843893
$ws.Range("G15").Value = -18.07964102843893
$ws.Range("H15").Value = -18.07964102843893
$ws.Range("I15").Value = -18.07964102843893
$ws.Range("J15").Value = -18.07964102843893
$ws.Range("K15").Value = -18.07964102843893
$ws.Range("B16").Value = -18.07964102843893
$ws.Range("C16").Value = -18.07964102843893
$ws.Range("D16").Value = -18.07964102843893
$ws.Range("E16").Value = -18.07964102843893
$ws.Range("F16").Value = -18.07964102843893
$ws.Range("G16").Value = -18.07964102843893
$ws.Range("H16").Value = -18.07964102843893
$ws.Range("I16").Value = -18.07964102843893
$ws.Range("J16").Value = -18.07964102843893
$ws.Range("K16").Value = -18.07964102843893
$ws.Range("B17").Value = -18.07964102843893
$ws.Range("C17").Value = 2.181798561022509
$ws.Range("D17").Value = 1.902646540422243
$ws.Range("E17").Value = -18.07964102843893
$ws.Range("F17").Value = -18.07964102843893
$ws.Range("G17").Value = -18.07964102843893
$ws.Range("H17").Value = 2.137958193987108
$ws.Range("I17").Value = 2.248855312606188
$ws.Range("J17").Value = -18.07964102843893
$ws.Range("K17").Value = -18.07964102843893
$ws.Range("B18").Value = -18.07964102843893
$ws.Range("C18").Value = -18.07964102843893
$ws.Range("D18").Value = -18.07964102843893
$ws.Range("E18").Value = -18.07964102843893
$ws.Range("F18").Value = -18.07964102843893
$ws.Range("G18").Value = -18.07964102843893
$ws.Range("H18").Value = 1.869679877722456
$ws.Range("I18").Value = 2.256439426059156
$ws.Range("J18").Value = -18.07964102843893
$ws.Range("K18").Value = -18.07964102843893
$ws.Range("B19").Value = -18.07964102843893
$ws.Range("C19").Value = -18.07964102843893
$ws.Range("D19").Value = 2.037262569201603
$ws.Range("E19").Value = -18.07964102843893
$ws.Range("F19").Value = -18.07964102843893
$ws.Range("G19").Value = -18.07964102843893
$ws.Range("H19").Value = 1.469190285075115
$ws.Range("I19").Value = 1.420292424611643
$ws.Range("J19").Value = -18.07964102843893
$ws.Range("K19").Value = -18.07964102843893
$ws.Range("B20").Value = -18.07964102843893
$ws.Range("C20").Value = 0.9557836696818689
$ws.Range("D20").Value = 1.530289759163075
$ws.Range("E20").Value = -18.07964102843893
$ws.Range("F20").Value = 3.38851811793336
$ws.Range("G20").Value = -18.07964102843893
$ws.Range("H20").Value = 1.887665454207403
$ws.Range("I20").Value = 1.580398401092949
$ws.Range("J20").Value = -18.07964102843893
$ws.Range("K20").Value = 2.464139420153086
$ws.Range("B21").Value = -18.07964102843893
$ws.Range("C21").Value = 1.341143675255374
$ws.Range("D21").Value = -18.07964102843893
$ws.Range("E21").Value = 1.678577305626386
$ws.Range("F21").Value = -18.07964102843893
$ws.Range("G21").Value = 2.858104182971637
$ws.Range("H21").Value = 1.567348280031454
$ws.Range("I21").Value = -18.07964102843893
$ws.Range("J21").Value = -18.07964102843893
$ws.Range("K21").Value = -18.07964102843893
